$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the attendance rows for Dipesh and Rubal (rows 4 and 5),
# leaving the formatting intact but clearing the cell values and
# their associated hyperlinks. Only the hyperlinks anchored in C4/C5
# should go away, so pick them out individually (re-querying the live
# collection for each one, since the collection re-indexes after every
# delete) rather than nuking the whole Hyperlinks collection.
$addrsToUnlink = @('$C$4', '$C$5')
foreach ($targetAddr in $addrsToUnlink) {
    $found = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $targetAddr) {
            $found = $h
            break
        }
    }
    if ($found -ne $null) {
        $found.Delete()
    }
}

$ws.Range("A4:C5").ClearContents()

# Update the active selection to G3.
$ws.Range("G3").Select()
